# Updated cryptos list (GitHub Actions data refresh): update Price (D)
# and Volume(1h) (E) columns for each coin row, and swap the
# EthereumClassic / Hedera rows (33 <-> 34).
#
# For cells whose new text would otherwise be auto-parsed by Excel as a
# number (losing formatting such as trailing zeros), the cell is forced
# to Text format first so the literal string is stored verbatim - exactly
# what typing that value into Excel and pre-formatting the cell as Text
# would produce.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.289.95"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.936.73"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.36"
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.90"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "2.932.66"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.75"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.46"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "65.343.85"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "3.424.53"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.03"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "2.933.54"
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.53"
$ws.Range("E20").Value = "  +11.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.45"
$ws.Range("E21").Value = "  -4.08%  "
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.34"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.10"
$ws.Range("E26").Value = "  -3.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  -5.67%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "0.0₃0997"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.10"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.971"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.53"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.14"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("E40").Value = "  -8.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.302"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("E43").Value = "  -7.68%  "
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "381.48"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "2.691.96"
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.01"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").Value = "  +4.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.106"
$ws.Range("E51").Value = "  -0.06%  "
